$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (46) of data to the Discounts sheet, following the pattern
# of the existing rows.
$row = 46
$ws.Cells.Item($row, 1).Value = 45
$ws.Cells.Item($row, 2).Value = 2
$ws.Cells.Item($row, 3).Value = 47
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = "System"

$ws.Cells.Item($row, 6).Value = "2025-03-04 07:04:08"

$ws.Cells.Item($row, 7).Value = 0
